$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2 (the "赣州·宅舞联萌·随舞动漫派对（免费活动)" entry) was removed from the
    # source data, so every following row moves up by one position.
    # (xlShiftUp = -4162). Note: in this runtime, Range.Delete shifts the whole
    # row regardless of the column span supplied, so column A's sequence
    # numbers get shifted too - they are fixed back up below.
    $ws.Range("A2:I2").Delete(-4162)

    # Column A holds a static running index (0,1,2,...) independent of the
    # row's data, so restore it to a clean 1..33 sequence for the 33
    # remaining data rows (rows 2..34).
    for ($i = 2; $i -le 34; $i++) {
        $ws.Cells.Item($i, 1).Value = $i - 1
    }
}
